# The "results" sheet had a "num_credito" row (row 2) that is removed.
# Deleting the row shifts every row below it up by one, which matches the
# target layout (dimension shrinks from A1:B9 to A1:B8) while keeping all
# remaining label/value pairs and their styling intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")
$ws.Rows.Item(2).Delete()
